$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextValue "D2" "91.902.87"
Set-TextValue "E2" "  -6.61%  "

Set-TextValue "D3" "3.298.53"
Set-TextValue "E3" "  -5.46%  "

Set-TextValue "E4" "  -0.05%  "

Set-TextValue "D5" "224.57"
Set-TextValue "E5" "  -11.20%  "

Set-TextValue "D6" "613.60"
Set-TextValue "E6" "  -7.99%  "

Set-TextValue "E7" "  -12.17%  "

Set-TextValue "D8" "0.369"
Set-TextValue "E8" "  -13.66%  "

Set-TextValue "E9" "  +0.10%  "

Set-TextValue "D10" "0.885"
Set-TextValue "E10" "  -16.59%  "

Set-TextValue "D11" "3.297.24"
Set-TextValue "E11" "  -5.46%  "

Set-TextValue "D12" "0.188"
Set-TextValue "E12" "  -11.13%  "

Set-TextValue "D13" "38.86"
Set-TextValue "E13" "  -14.67%  "

Set-TextValue "D14" "91.699.16"
Set-TextValue "E14" "  -6.64%  "

Set-TextValue "D15" "5.72"
Set-TextValue "E15" "  -8.58%  "

Set-TextValue "D16" "3.909.35"
Set-TextValue "E16" "  -5.70%  "

Set-TextValue "D17" "0.0000235"
Set-TextValue "E17" "  -10.01%  "

Set-TextValue "D18" "3.294.17"
Set-TextValue "E18" "  -5.53%  "

Set-TextValue "D19" "7.61"
Set-TextValue "E19" "  -15.38%  "

Set-TextValue "D20" "16.30"
Set-TextValue "E20" "  -13.02%  "

Set-TextValue "D21" "10.56"
Set-TextValue "E21" "  -10.11%  "

Set-TextValue "D22" "474.84"
Set-TextValue "E22" "  -9.19%  "

Set-TextValue "D23" "0.427"
Set-TextValue "E23" "  -19.15%  "

Set-TextValue "D24" "3.00"
Set-TextValue "E24" "  -12.42%  "

Set-TextValue "D25" "0.0000175"
Set-TextValue "E25" "  -14.02%  "

Set-TextValue "D26" "5.94"
Set-TextValue "E26" "  -12.64%  "

Set-TextValue "D27" "87.79"
Set-TextValue "E27" "  -10.85%  "

Set-TextValue "B28" "Aptos"
Set-TextValue "C28" "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
Set-TextValue "D28" "11.04"
Set-TextValue "E28" "  -13.41%  "

Set-TextValue "B29" "Dai"
Set-TextValue "C29" "https://coinranking.com/coin/MoTuySvg7+dai-dai"
Set-TextValue "D29" "1.00"
Set-TextValue "E29" "  +0.01%  "

Set-TextValue "B30" "InternetComputer(DFINITY)"
Set-TextValue "C30" "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
Set-TextValue "D30" "10.88"
Set-TextValue "E30" "  -12.75%  "

Set-TextValue "B31" "Binance-PegBSC-USD"
Set-TextValue "C31" "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
Set-TextValue "D31" "1.01"
Set-TextValue "E31" "  +0.48%  "

Set-TextValue "D32" "2.53"
Set-TextValue "E32" "  -11.77%  "

Set-TextValue "B33" "Hedera"
Set-TextValue "C33" "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextValue "D33" "0.127"
Set-TextValue "E33" "  -13.15%  "

Set-TextValue "B34" "Cronos"
Set-TextValue "C34" "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
Set-TextValue "D34" "0.165"
Set-TextValue "E34" "  -13.28%  "

Set-TextValue "B35" "EthereumClassic"
Set-TextValue "C35" "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
Set-TextValue "D35" "27.63"
Set-TextValue "E35" "  -11.34%  "

Set-TextValue "B36" "PolygonEcosystemToken"
Set-TextValue "C36" "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
Set-TextValue "D36" "0.502"
Set-TextValue "E36" "  -15.98%  "

Set-TextValue "B37" "USDe"
Set-TextValue "C37" "https://coinranking.com/coin/exbfr2U-0+usde-usde"
Set-TextValue "D37" "1.00"
Set-TextValue "E37" "  -0.06%  "

Set-TextValue "B38" "Bittensor"
Set-TextValue "C38" "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
Set-TextValue "D38" "508.60"
Set-TextValue "E38" "  -3.35%  "

Set-TextValue "B39" "RenderToken"
Set-TextValue "C39" "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
Set-TextValue "D39" "7.15"
Set-TextValue "E39" "  -10.03%  "

Set-TextValue "B40" "Fetch.AI"
Set-TextValue "C40" "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
Set-TextValue "D40" "1.34"
Set-TextValue "E40" "  -11.56%  "

Set-TextValue "B41" "Kaspa"
Set-TextValue "C41" "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
Set-TextValue "D41" "0.142"
Set-TextValue "E41" "  -8.84%  "

Set-TextValue "B42" "WhiteBITCoin"
Set-TextValue "C42" "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
Set-TextValue "D42" "23.96"
Set-TextValue "E42" "  -1.84%  "

Set-TextValue "D43" "0.845"
Set-TextValue "E43" "  -7.61%  "

Set-TextValue "B44" "MantraDAO"
Set-TextValue "C44" "https://coinranking.com/coin/cTdD8lD-6+mantradao-om"
Set-TextValue "D44" "3.48"
Set-TextValue "E44" "  -4.36%  "

Set-TextValue "D45" "1.61"
Set-TextValue "E45" "  -8.68%  "

Set-TextValue "B46" "Filecoin"
Set-TextValue "C46" "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextValue "D46" "5.23"
Set-TextValue "E46" "  -9.75%  "

Set-TextValue "B47" "Stacks"
Set-TextValue "C47" "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
Set-TextValue "D47" "2.07"
Set-TextValue "E47" "  -7.06%  "

Set-TextValue "B48" "OKB"
Set-TextValue "C48" "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
Set-TextValue "D48" "51.33"
Set-TextValue "E48" "  -7.45%  "

Set-TextValue "B49" "VeChain"
Set-TextValue "C49" "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextValue "D49" "0.0379"
Set-TextValue "E49" "  -12.94%  "

Set-TextValue "B50" "Cosmos"
Set-TextValue "C50" "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
Set-TextValue "D50" "7.57"
Set-TextValue "E50" "  -13.29%  "

Set-TextValue "D51" "2.95"
Set-TextValue "E51" "  -8.86%  "
